# Generate Report for Handback
# Updates status/timestamps on the Overview, zh-cn, and de-de sheets to
# reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
# Rows 2 and 5 both reference the same timestamp value.
$wsOverview.Range("G2").Value = "2016-08-13 06:19:20"
$wsOverview.Range("G5").Value = "2016-08-13 06:19:20"

# --- zh-cn sheet ---
# "Priority" column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# "Correspond Handoff Datetime" column (H)
$wsZhCn.Range("H2").Value = "2016-08-13 06:19:12"
$wsZhCn.Range("H5").Value = "2016-08-13 06:19:12"

# "Correspond Handback DateTime" column (K)
$wsZhCn.Range("K2").Value = "2016-08-13 06:19:41"
$wsZhCn.Range("K5").Value = "2016-08-13 06:19:41"

# --- de-de sheet ---
# "Priority" column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# "Correspond Handoff Datetime" column (H)
$wsDeDe.Range("H2").Value = "2016-08-13 06:19:20"
$wsDeDe.Range("H5").Value = "2016-08-13 06:19:20"

# "Correspond Handback DateTime" column (K)
$wsDeDe.Range("K2").Value = "2016-08-13 06:19:50"
$wsDeDe.Range("K5").Value = "2016-08-13 06:19:50"
